$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 124, shifting the existing
# rows 124-133 down to 126-135 (dates/values stay attached to their rows).
$ws.Rows.Item(124).Resize(2).Insert()

# --- New row 124: Alcachofa / Española ---
$ws.Cells.Item(124,1).Value  = 7
$ws.Cells.Item(124,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(124,3).Value  = "Ñuble"
$ws.Cells.Item(124,4).Value  = 45194
$ws.Cells.Item(124,5).Value  = 16
$ws.Cells.Item(124,6).Value  = 100112013
$ws.Cells.Item(124,7).Value  = "Alcachofa"
$ws.Cells.Item(124,8).Value  = "Española"
$ws.Cells.Item(124,9).Value  = "Primera"
$ws.Cells.Item(124,10).Value = 50
$ws.Cells.Item(124,11).Value = 12000
$ws.Cells.Item(124,12).Value = 12000
$ws.Cells.Item(124,13).Value = 12000
$ws.Cells.Item(124,14).Value = "$/caja 30 unidades"
$ws.Cells.Item(124,15).Value = "Provincia de Limarí"
$ws.Cells.Item(124,16).Value = 400
$ws.Cells.Item(124,17).Value = 30
$ws.Cells.Item(124,18).Value = "Hortaliza"

# --- New row 125: Alcachofa / Madrigal ---
$ws.Cells.Item(125,1).Value  = 7
$ws.Cells.Item(125,2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(125,3).Value  = "Ñuble"
$ws.Cells.Item(125,4).Value  = 45194
$ws.Cells.Item(125,5).Value  = 16
$ws.Cells.Item(125,6).Value  = 100112013
$ws.Cells.Item(125,7).Value  = "Alcachofa"
$ws.Cells.Item(125,8).Value  = "Madrigal"
$ws.Cells.Item(125,9).Value  = "Primera"
$ws.Cells.Item(125,10).Value = 50
$ws.Cells.Item(125,11).Value = 12000
$ws.Cells.Item(125,12).Value = 12000
$ws.Cells.Item(125,13).Value = 12000
$ws.Cells.Item(125,14).Value = "$/caja 40 unidades"
$ws.Cells.Item(125,15).Value = "Provincia del Elquí"
$ws.Cells.Item(125,16).Value = 300
$ws.Cells.Item(125,17).Value = 40
$ws.Cells.Item(125,18).Value = "Hortaliza"
